# Apply KHL referees stats update (khl_referees_stats_1369.xlsx)
#
# Two sheets - "Главные" (2nd tab) and "Линейные" (3rd tab) - both keep a
# per-referee PIM (penalty-in-minutes) stats table in columns A:AA. This
# edit refreshes a subset of rows with newer cumulative totals (more games
# played / more penalty minutes since the last snapshot) and stamps every
# data row (2-26) with the new as_of_utc pull time in column AA.

$wb = $excel.ActiveWorkbook

$wsGlavnye = $wb.Worksheets.Item(2)
$wsLineynye = $wb.Worksheets.Item(3)

$glavnyeUpdates = @{
    "2" = @{ "AA" = "2025-11-08 03:03:59"; "C" = 22; "D" = 510; "E" = 222; "F" = 288; "G" = 23.18; "H" = 10.09; "I" = 13.09; "J" = 96; "K" = 114 }
    "3" = @{ "AA" = "2025-11-08 03:03:59"; "C" = 22; "D" = 401; "E" = 175; "F" = 226; "G" = 18.23; "I" = 10.27; "J" = 85; "K" = 93 }
    "4" = @{ "AA" = "2025-11-08 03:03:59" }
    "5" = @{ "AA" = "2025-11-08 03:03:59" }
    "6" = @{ "AA" = "2025-11-08 03:03:59"; "C" = 22; "D" = 405; "E" = 169; "F" = 236; "G" = 18.41; "H" = 7.68; "I" = 10.73; "J" = 77; "K" = 98; "W" = 14 }
    "7" = @{ "AA" = "2025-11-08 03:03:59" }
    "8" = @{ "AA" = "2025-11-08 03:03:59" }
    "9" = @{ "AA" = "2025-11-08 03:03:59" }
    "10" = @{ "AA" = "2025-11-08 03:03:59" }
    "11" = @{ "AA" = "2025-11-08 03:03:59" }
    "12" = @{ "AA" = "2025-11-08 03:03:59"; "C" = 14; "D" = 247; "E" = 102; "F" = 145; "G" = 17.64; "H" = 7.29; "I" = 10.36; "J" = 41; "K" = 50; "X" = 9 }
    "13" = @{ "AA" = "2025-11-08 03:03:59" }
    "14" = @{ "AA" = "2025-11-08 03:03:59" }
    "15" = @{ "AA" = "2025-11-08 03:03:59" }
    "16" = @{ "AA" = "2025-11-08 03:03:59" }
    "17" = @{ "AA" = "2025-11-08 03:03:59" }
    "18" = @{ "AA" = "2025-11-08 03:03:59" }
    "19" = @{ "AA" = "2025-11-08 03:03:59"; "C" = 17; "D" = 312; "E" = 144; "F" = 168; "G" = 18.35; "H" = 8.470000000000001; "I" = 9.880000000000001; "J" = 67; "K" = 69; "X" = 4 }
    "20" = @{ "AA" = "2025-11-08 03:03:59" }
    "21" = @{ "AA" = "2025-11-08 03:03:59" }
    "22" = @{ "AA" = "2025-11-08 03:03:59"; "C" = 16; "D" = 314; "E" = 124; "F" = 190; "G" = 19.63; "H" = 7.75; "I" = 11.88; "J" = 62; "K" = 65 }
    "23" = @{ "AA" = "2025-11-08 03:03:59" }
    "24" = @{ "AA" = "2025-11-08 03:03:59"; "C" = 22; "D" = 344; "E" = 165; "F" = 179; "G" = 15.64; "H" = 7.5; "I" = 8.140000000000001; "J" = 80; "K" = 87 }
    "25" = @{ "AA" = "2025-11-08 03:03:59"; "C" = 22; "D" = 384; "E" = 194; "F" = 190; "G" = 17.45; "H" = 8.82; "I" = 8.640000000000001; "J" = 92; "K" = 90; "W" = 8 }
    "26" = @{ "AA" = "2025-11-08 03:03:59" }
}

$lineynyeUpdates = @{
    "2" = @{ "AA" = "2025-11-08 03:03:59" }
    "3" = @{ "AA" = "2025-11-08 03:03:59" }
    "4" = @{ "AA" = "2025-11-08 03:03:59" }
    "5" = @{ "AA" = "2025-11-08 03:03:59" }
    "6" = @{ "AA" = "2025-11-08 03:03:59" }
    "7" = @{ "AA" = "2025-11-08 03:03:59" }
    "8" = @{ "AA" = "2025-11-08 03:03:59"; "C" = 20; "D" = 292; "E" = 120; "F" = 172; "G" = 14.6; "I" = 8.6; "J" = 55; "K" = 71 }
    "9" = @{ "AA" = "2025-11-08 03:03:59" }
    "10" = @{ "AA" = "2025-11-08 03:03:59" }
    "11" = @{ "AA" = "2025-11-08 03:03:59"; "C" = 13; "D" = 191; "E" = 92; "F" = 99; "G" = 14.69; "H" = 7.08; "I" = 7.62; "J" = 46; "K" = 47 }
    "12" = @{ "AA" = "2025-11-08 03:03:59" }
    "13" = @{ "AA" = "2025-11-08 03:03:59" }
    "14" = @{ "AA" = "2025-11-08 03:03:59" }
    "15" = @{ "AA" = "2025-11-08 03:03:59"; "C" = 19; "D" = 391; "E" = 205; "F" = 186; "G" = 20.58; "H" = 10.79; "I" = 9.789999999999999; "J" = 80; "K" = 73; "X" = 9 }
    "16" = @{ "AA" = "2025-11-08 03:03:59"; "C" = 21; "D" = 360; "E" = 174; "F" = 186; "G" = 17.14; "H" = 8.289999999999999; "I" = 8.859999999999999; "J" = 82; "K" = 88 }
    "17" = @{ "AA" = "2025-11-08 03:03:59"; "C" = 10; "D" = 144; "E" = 86; "F" = 58; "G" = 14.4; "H" = 8.6; "I" = 5.8; "J" = 43; "K" = 29 }
    "18" = @{ "AA" = "2025-11-08 03:03:59"; "C" = 23; "D" = 403; "E" = 192; "F" = 211; "G" = 17.52; "H" = 8.35; "I" = 9.17; "J" = 91; "K" = 88; "W" = 16 }
    "19" = @{ "AA" = "2025-11-08 03:03:59" }
    "20" = @{ "AA" = "2025-11-08 03:03:59" }
    "21" = @{ "AA" = "2025-11-08 03:03:59" }
    "22" = @{ "AA" = "2025-11-08 03:03:59" }
    "23" = @{ "AA" = "2025-11-08 03:03:59"; "C" = 13; "D" = 201; "E" = 94; "F" = 107; "G" = 15.46; "H" = 7.23; "I" = 8.23; "J" = 47; "K" = 51; "X" = 4 }
    "24" = @{ "AA" = "2025-11-08 03:03:59"; "C" = 23; "D" = 416; "E" = 167; "F" = 249; "G" = 18.09; "H" = 7.26; "I" = 10.83; "J" = 76; "K" = 97; "W" = 14 }
    "25" = @{ "AA" = "2025-11-08 03:03:59" }
    "26" = @{ "AA" = "2025-11-08 03:03:59" }
}

foreach ($rowKey in $glavnyeUpdates.Keys) {
    $cellData = $glavnyeUpdates[$rowKey]
    foreach ($colKey in $cellData.Keys) {
        $addr = "$colKey$rowKey"
        $wsGlavnye.Range($addr).Value = $cellData[$colKey]
    }
}

foreach ($rowKey in $lineynyeUpdates.Keys) {
    $cellData = $lineynyeUpdates[$rowKey]
    foreach ($colKey in $cellData.Keys) {
        $addr = "$colKey$rowKey"
        $wsLineynye.Range($addr).Value = $cellData[$colKey]
    }
}

Write-Host "KHL referees stats update applied."
